# Apply the "6.0.0" FHIR IG republish update to the StructureDefinition
# exchange-plan-indicator workbook:
#  - Metadata sheet: bump Version/Date, fill in Publisher, replace the
#    duplicated "Contact" row with a "Jurisdiction" row, and drop the
#    leftover duplicate row underneath it.
#  - Elements sheet: give the root Extension row its real Short/Definition
#    text instead of the generic placeholder.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date bumped to the new publish timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank; now populated
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 changes from Contact/"No display for ContactDetail" to Jurisdiction
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate Contact row with no real content; remove it so
# everything below shifts up by one (A1:B21 -> A1:B20)
$meta.Rows.Item(11).Delete()

$elements = $wb.Worksheets.Item("Elements")

# The root Extension row (row 2) gets its real Short/Definition text
$elements.Range("K2").Value = "Exchange Plan Indicator"
$elements.Range("L2").Value = "Indicator that identifies if the plan is available on a healthcare exchange marketplace"
